$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 becomes the former row 3's "BusinessMan" record, with a couple of
# fields (name and education) edited to new values.
$ws.Range("A1").Value = "AbdullahAlMamun"
$ws.Range("B1").Value = "BusinessMan"
$ws.Range("C1").Value = 85000
$ws.Range("D1").Value = "MA"
$ws.Range("E1").Value = "Chittagong"
$ws.Range("F1").Value = 25.3698
$ws.Range("G1").Value = 15000

# Row 2 (the "Driver" record) is removed entirely.
$ws.Range("A2:G2").ClearContents() | Out-Null

# Row 3's data is cleared out too, leaving only its existing formatted
# (wrap-text) cells behind in D3/E3.
$ws.Range("A3:G3").ClearContents() | Out-Null

# Widen the two new descriptive columns. The host snaps ColumnWidth to a
# 1/6-character grid on save (stored = round((input + 5/6) * 6) / 6), so
# these inputs are the bucket centers that land the saved width as close
# as the host's precision allows to the target 18.5546875 / 13.21875.
$ws.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws.Columns.Item(2).ColumnWidth = 12.333333333333332

# Move the active cell selection.
$ws.Range("F7").Select() | Out-Null
